$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 586.6316
$ws.Range("I2").Value = 282
$ws.Range("J2").Value = 1439.6
$ws.Range("K2").Value = 282
$ws.Range("L2").Value = 1439.6
$ws.Range("M2").Value = -169
$ws.Range("N2").Value = -1665.6
$ws.Range("H9").Value = 319.15384
$ws.Range("I9").Value = 309.9
$ws.Range("K9").Value = 309.9
$ws.Range("M9").Value = -140.9
$ws.Range("H12").Value = 20810.875
$ws.Range("I12").Value = 16099.8
$ws.Range("K12").Value = 16099.8
$ws.Range("M12").Value = -15929.8
$ws.Range("H43").Value = 1697.3
$ws.Range("I43").Value = 1715.1428
$ws.Range("J43").Value = 1655.6666
$ws.Range("K43").Value = 1715.1428
$ws.Range("L43").Value = 1655.6666
$ws.Range("M43").Value = -1646.1428
$ws.Range("N43").Value = -1793.6666
$ws.Range("H59").Value = 72.5
$ws.Range("I59").Value = 72.5
$ws.Range("K59").Value = 217.5
$ws.Range("M59").Value = 339.5
$ws.Range("H80").Value = 4141.952
$ws.Range("I80").Value = 1167.75
$ws.Range("J80").Value = 5972.231
$ws.Range("K80").Value = 3503.25
$ws.Range("L80").Value = 17916.693
$ws.Range("M80").Value = -2505.25
$ws.Range("N80").Value = -19912.693
$ws.Range("H83").Value = 4141.952
$ws.Range("I83").Value = 1167.75
$ws.Range("J83").Value = 5972.231
$ws.Range("K83").Value = 10509.75
$ws.Range("L83").Value = 53750.079
$ws.Range("M83").Value = -5517.75
$ws.Range("N83").Value = -63734.079
$ws.Range("H86").Value = 2460
$ws.Range("J86").Value = 2669.9
$ws.Range("L86").Value = 2669.9
$ws.Range("N86").Value = -4915.9
$ws.Range("H89").Value = 2460
$ws.Range("J89").Value = 2669.9
$ws.Range("L89").Value = 13349.5
$ws.Range("N89").Value = -24581.5
$ws.Range("H140").Value = 60780
$ws.Range("J140").Value = 60780
$ws.Range("L140").Value = 60780
$ws.Range("N140").Value = -71140
$ws.Range("H141").Value = 4816.3447
$ws.Range("I141").Value = 5255.522
$ws.Range("K141").Value = 15766.566
$ws.Range("M141").Value = -10586.566

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 87959.96000000001
$ws.Range("I74").Value = 57673.953
$ws.Range("K74").Value = 57673.953
$ws.Range("M74").Value = -56799.953
$ws.Range("H77").Value = 87959.96000000001
$ws.Range("I77").Value = 57673.953
$ws.Range("K77").Value = 288369.765
$ws.Range("M77").Value = -284001.765
$ws.Range("H88").Value = 1061.6
$ws.Range("J88").Value = 998
$ws.Range("L88").Value = 998
$ws.Range("N88").Value = -1810
$ws.Range("H91").Value = 1061.6
$ws.Range("J91").Value = 998
$ws.Range("L91").Value = 998
$ws.Range("N91").Value = -3806

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1466.8334
$ws.Range("I20").Value = 1246.2
$ws.Range("J20").Value = 1908.1
$ws.Range("K20").Value = 1246.2
$ws.Range("L20").Value = 1908.1
$ws.Range("M20").Value = -999.2
$ws.Range("N20").Value = -2402.1
$ws.Range("H47").Value = 310000
$ws.Range("J47").Value = 310000
$ws.Range("L47").Value = 310000
$ws.Range("N47").Value = -311040
$ws.Range("H86").Value = 8635.677
$ws.Range("I86").Value = 7321.84
$ws.Range("J86").Value = 12285.223
$ws.Range("K86").Value = 7321.84
$ws.Range("L86").Value = 12285.223
$ws.Range("M86").Value = -6198.84
$ws.Range("N86").Value = -14531.223
$ws.Range("H89").Value = 8635.677
$ws.Range("I89").Value = 7321.84
$ws.Range("J89").Value = 12285.223
$ws.Range("K89").Value = 36609.2
$ws.Range("L89").Value = 61426.115
$ws.Range("M89").Value = -30993.2
$ws.Range("N89").Value = -72658.11499999999
$ws.Range("H105").Value = 2193.4285
$ws.Range("I105").Value = 2193.4285
$ws.Range("K105").Value = 2193.4285
$ws.Range("M105").Value = -446.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3153.3333
$ws.Range("I31").Value = 1954.7307
$ws.Range("K31").Value = 1954.7307
$ws.Range("M31").Value = -1659.7307
$ws.Range("H34").Value = 3153.3333
$ws.Range("I34").Value = 1954.7307
$ws.Range("K34").Value = 1954.7307
$ws.Range("M34").Value = -1752.7307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 13137.375
$ws.Range("I11").Value = 10050
$ws.Range("J11").Value = 14166.5
$ws.Range("K11").Value = 30150
$ws.Range("L11").Value = 42499.5
$ws.Range("M11").Value = -30010
$ws.Range("N11").Value = -42779.5
$ws.Range("H26").Value = 110
$ws.Range("I26").Value = 110
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 330
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -42
$ws.Range("N26").Value = ""
$ws.Range("H37").Value = 41410.445
$ws.Range("J37").Value = 41410.445
$ws.Range("L37").Value = 124231.335
$ws.Range("N37").Value = -124455.335
$ws.Range("H54").Value = 2330.6667
$ws.Range("J54").Value = 2330.6667
$ws.Range("L54").Value = 6992.000100000001
$ws.Range("N54").Value = -8110.000100000001
$ws.Range("H56").Value = 10006325
$ws.Range("I56").Value = 10006325
$ws.Range("K56").Value = 10006325
$ws.Range("M56").Value = -10005795
$ws.Range("H62").Value = 39014
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 39014
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 117042
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -118414
$ws.Range("H65").Value = 39014
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 39014
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 351126
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -357990
$ws.Range("H126").Value = 830
$ws.Range("I126").Value = 830
$ws.Range("K126").Value = 2490
$ws.Range("M126").Value = 2450
$ws.Range("H129").Value = 1443.6666
$ws.Range("J129").Value = 999.5
$ws.Range("L129").Value = 2998.5
$ws.Range("N129").Value = -12998.5
$ws.Range("H131").Value = 13023213
$ws.Range("I131").Value = 11905595
$ws.Range("K131").Value = 35716785
$ws.Range("M131").Value = -35711745

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 51.625
$ws.Range("I2").Value = 47.6
$ws.Range("K2").Value = 47.6
$ws.Range("M2").Value = 65.40000000000001
$ws.Range("H11").Value = 10000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = ""
$ws.Range("H70").Value = 6500
$ws.Range("I70").Value = 6500
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6500
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -6230
$ws.Range("N70").Value = ""
$ws.Range("H73").Value = 6500
$ws.Range("I73").Value = 6500
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6500
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -5564
$ws.Range("N73").Value = ""
$ws.Range("H122").Value = 76186.48
$ws.Range("I122").Value = 117353.69
$ws.Range("K122").Value = 352061.07
$ws.Range("M122").Value = -349611.07

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 9000
$ws.Range("I20").Value = 9000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -8774
$ws.Range("N20").Value = ""
$ws.Range("H22").Value = 1853.5385
$ws.Range("I22").Value = 1249.75
$ws.Range("J22").Value = 2121.889
$ws.Range("K22").Value = 1249.75
$ws.Range("L22").Value = 2121.889
$ws.Range("M22").Value = -954.75
$ws.Range("N22").Value = -2711.889
$ws.Range("H27").Value = 1853.5385
$ws.Range("I27").Value = 1249.75
$ws.Range("J27").Value = 2121.889
$ws.Range("K27").Value = 1249.75
$ws.Range("L27").Value = 2121.889
$ws.Range("M27").Value = -1142.75
$ws.Range("N27").Value = -2335.889
$ws.Range("H40").Value = 6452.0713
$ws.Range("J40").Value = 11107.333
$ws.Range("L40").Value = 11107.333
$ws.Range("N40").Value = -11379.333
$ws.Range("H68").Value = 3500
$ws.Range("I68").Value = 4000
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 4000
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -3251
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 3500
$ws.Range("I71").Value = 4000
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 20000
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -16256
$ws.Range("N71").Value = -22488
$ws.Range("H132").Value = 9551.565000000001
$ws.Range("I132").Value = 11512.2
$ws.Range("J132").Value = 5875.375
$ws.Range("K132").Value = 34536.60000000001
$ws.Range("L132").Value = 17626.125
$ws.Range("M132").Value = -32006.60000000001
$ws.Range("N132").Value = -22686.125
$ws.Range("H136").Value = 62354.5
$ws.Range("I136").Value = 85377.25
$ws.Range("K136").Value = 256131.75
$ws.Range("M136").Value = -253581.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2351.25
$ws.Range("I100").Value = 2351.25
$ws.Range("K100").Value = 4702.5
$ws.Range("M100").Value = -4161.5
$ws.Range("H109").Value = 49998.5
$ws.Range("J109").Value = 49998.5
$ws.Range("L109").Value = 49998.5
$ws.Range("N109").Value = -52772.5
$ws.Range("H132").Value = 319729.25
$ws.Range("I132").Value = 8703.571
$ws.Range("J132").Value = 913505.5600000001
$ws.Range("K132").Value = 26110.713
$ws.Range("L132").Value = 2740516.68
$ws.Range("M132").Value = -23580.713
$ws.Range("N132").Value = -2745576.68
$ws.Range("H136").Value = 2549.3
$ws.Range("I136").Value = 1201.6818
$ws.Range("K136").Value = 3605.0454
$ws.Range("M136").Value = -1055.0454
